$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing data row (row 2), pushing the
# current "Oogong Elementary School" row down to row 3.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the "Gatid Elementary School" data.
$ws.Range("A2").Value = "Gatid Elementary School"
$ws.Range("B2").Value = "Laguna"
$ws.Range("C2").Value = "Santa Cruz"
$ws.Range("D2").Value = 123461
$ws.Range("E2").Value = "Gatid, Santa Cruz, Laguna"
$ws.Range("F2").Value = "Barangay Gatid"
$ws.Range("G2").Value = "example name"

# The contact number looks like a number (leading zero), so assign it via a
# text formula and then convert the formula to a static value in place -
# this keeps it stored as a shared string (matching the source data) without
# introducing a new cell style (no "@" text number format needed).
$ws.Range("H2").Formula = "=""09123345353"""
$ws.Range("H2").Copy()
$ws.Range("H2").PasteSpecial(-4163)

$ws.Range("I2").Value = 4
